$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet
$ws.Name = "DanhSachChucVu"

# Header row
$ws.Range("A1").Value = "Mã chức vụ"
$ws.Range("B1").Value = "Tên chức vụ"
$ws.Range("C1").Value = "Trạng thái hiển thị"

# Data rows
$ws.Range("A2").Value = "BCS"
$ws.Range("B2").Value = "Ban cán sự"
$ws.Range("C2").Value = "Hiển thị"

$ws.Range("A3").Value = "CV"
$ws.Range("B3").Value = "Cố vấn học tập "
$ws.Range("C3").Value = "Hiển thị"

$ws.Range("A4").Value = "QL"
$ws.Range("B4").Value = "Quản lý"
$ws.Range("C4").Value = "Hiển thị"

$ws.Range("A5").Value = "SV"
$ws.Range("B5").Value = "Sinh viên"
$ws.Range("C5").Value = "Hiển thị"
